# ProjectTest65/Excel/Login.xlsx - chore: clean up code structure and remove redundant code blocks
#
# The underlying test-results sheet ("Login") gets a handful of data-entry
# updates: a previously-blank result row is filled in, and a few rows' Pass/Fail
# outcome (and accompanying actual-result message) are corrected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")

# Row 2 ("Execute" = Y, login succeeded test case) was missing its
# ActualResult/Result columns - fill them in.
$ws.Range("F2").Value = "No message found"
$ws.Range("G2").Value = "Fail"

# Mark these test cases as executed (N -> Y).
$ws.Range("A3").Value = "Y"
$ws.Range("A8").Value = "Y"
$ws.Range("A9").Value = "Y"

# Row 8: actual result / outcome corrected to reflect a failed login attempt.
$ws.Range("F8").Value = "Email หรือ Passwordไม่ถูกต้อง กรุณาลองใหม่อีกครั้ง"
$ws.Range("G8").Value = "Fail"

# Row 10: actual result / outcome corrected to reflect a failed login attempt.
$ws.Range("F10").Value = "No message found"
$ws.Range("G10").Value = "Fail"

# Leave the cursor where the author last left it.
$ws.Range("G11").Select()
